$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6125
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 6125
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6125
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7093

$ws.Range("H62").Value = 7713.0835
$ws.Range("I62").Value = 7426.3335
$ws.Range("K62").Value = 7426.3335
$ws.Range("M62").Value = -6802.3335

$ws.Range("H65").Value = 7713.0835
$ws.Range("I65").Value = 7426.3335
$ws.Range("K65").Value = 37131.6675
$ws.Range("M65").Value = -34011.6675

$ws.Range("H92").Value = 1225.3438
$ws.Range("I92").Value = 369.12
$ws.Range("K92").Value = 369.12
$ws.Range("M92").Value = 878.88

$ws.Range("H96").Value = 225.11765
$ws.Range("I96").Value = 264.81818
$ws.Range("J96").Value = 152.33333
$ws.Range("K96").Value = 794.45454
$ws.Range("L96").Value = 456.99999
$ws.Range("M96").Value = 578.54546
$ws.Range("N96").Value = -3202.99999

$ws.Range("H137").Value = 32300.758
$ws.Range("I137").Value = 39192.363
$ws.Range("K137").Value = 117577.089
$ws.Range("M137").Value = -115027.089

$ws.Range("H138").Value = 2804.9858
$ws.Range("I138").Value = 1287.0834
$ws.Range("J138").Value = 3580.0852
$ws.Range("K138").Value = 3861.2502
$ws.Range("L138").Value = 10740.2556
$ws.Range("M138").Value = 1278.7498
$ws.Range("N138").Value = -21020.2556

$ws.Range("H141").Value = 7382.727
$ws.Range("I141").Value = 7643.8096
$ws.Range("K141").Value = 22931.4288
$ws.Range("M141").Value = -17751.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1825668.8
$ws.Range("I2").Value = 3143092.2
$ws.Range("J2").Value = 1543.8462
$ws.Range("K2").Value = 3143092.2
$ws.Range("L2").Value = 1543.8462
$ws.Range("M2").Value = -3142979.2
$ws.Range("N2").Value = -1769.8462

$ws.Range("H32").Value = 9975.115
$ws.Range("I32").Value = 5503.3555
$ws.Range("J32").Value = 22551.938
$ws.Range("K32").Value = 5503.3555
$ws.Range("L32").Value = 22551.938
$ws.Range("M32").Value = -5216.3555
$ws.Range("N32").Value = -23125.938

$ws.Range("H45").Value = 6692720.5
$ws.Range("I45").Value = 9052129
$ws.Range("K45").Value = 9052129
$ws.Range("M45").Value = -9051752

$ws.Range("H61").Value = 4050
$ws.Range("I61").Value = 4005
$ws.Range("K61").Value = 4005
$ws.Range("M61").Value = -3793

$ws.Range("H102").Value = 13895228
$ws.Range("I102").Value = 27784448
$ws.Range("J102").Value = 6007.3335
$ws.Range("K102").Value = 27784448
$ws.Range("L102").Value = 6007.3335
$ws.Range("M102").Value = -27782826
$ws.Range("N102").Value = -9251.333500000001

$ws.Range("H116").Value = 1825668.8
$ws.Range("I116").Value = 3143092.2
$ws.Range("J116").Value = 1543.8462
$ws.Range("K116").Value = 3143092.2
$ws.Range("L116").Value = 1543.8462
$ws.Range("M116").Value = -3140798.2
$ws.Range("N116").Value = -6131.8462

$ws.Range("H136").Value = 4050
$ws.Range("I136").Value = 4005
$ws.Range("K136").Value = 12015
$ws.Range("M136").Value = -9465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1825668.8
$ws.Range("I3").Value = 3143092.2
$ws.Range("J3").Value = 1543.8462
$ws.Range("K3").Value = 3143092.2
$ws.Range("L3").Value = 1543.8462
$ws.Range("M3").Value = -3142978.2
$ws.Range("N3").Value = -1771.8462

$ws.Range("H24").Value = 2153.5
$ws.Range("J24").Value = 2199
$ws.Range("L24").Value = 2199
$ws.Range("N24").Value = -2669

$ws.Range("H86").Value = 4001123.8
$ws.Range("I86").Value = 6667672.5
$ws.Range("J86").Value = 1300.8
$ws.Range("K86").Value = 6667672.5
$ws.Range("L86").Value = 1300.8
$ws.Range("M86").Value = -6666549.5
$ws.Range("N86").Value = -3546.8

$ws.Range("H89").Value = 4001123.8
$ws.Range("I89").Value = 6667672.5
$ws.Range("J89").Value = 1300.8
$ws.Range("K89").Value = 33338362.5
$ws.Range("L89").Value = 6504
$ws.Range("M89").Value = -33332746.5
$ws.Range("N89").Value = -17736

$ws.Range("H94").Value = 4352666
$ws.Range("I94").Value = 7143555
$ws.Range("J94").Value = 11283.223
$ws.Range("K94").Value = 7143555
$ws.Range("L94").Value = 11283.223
$ws.Range("M94").Value = -7143104
$ws.Range("N94").Value = -12185.223

$ws.Range("H105").Value = 5209051.5
$ws.Range("I105").Value = 5682492.5
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 5682492.5
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -5680745.5
$ws.Range("N105").Value = -4694

$ws.Range("H107").Value = 2860054.2
$ws.Range("I107").Value = 3970408.8
$ws.Range("J107").Value = 4857
$ws.Range("K107").Value = 3970408.8
$ws.Range("L107").Value = 4857
$ws.Range("M107").Value = -3968488.8
$ws.Range("N107").Value = -8697

$ws.Range("H134").Value = 3154.7334
$ws.Range("I134").Value = 1371.3793
$ws.Range("K134").Value = 4114.1379
$ws.Range("M134").Value = -1579.1379

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 57333.332
$ws.Range("I57").Value = 57333.332
$ws.Range("K57").Value = 57333.332
$ws.Range("M57").Value = -56773.332

$ws.Range("H58").Value = 4462.241
$ws.Range("I58").Value = 5306
$ws.Range("K58").Value = 5306
$ws.Range("M58").Value = -5103

$ws.Range("H59").Value = 23534.666
$ws.Range("I59").Value = 15052
$ws.Range("K59").Value = 15052
$ws.Range("M59").Value = -13907

$ws.Range("H132").Value = 38537.453
$ws.Range("I132").Value = 24121.695
$ws.Range("K132").Value = 72365.08499999999
$ws.Range("M132").Value = -69835.08499999999

$ws.Range("H136").Value = 4462.241
$ws.Range("I136").Value = 5306
$ws.Range("K136").Value = 15918
$ws.Range("M136").Value = -13368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 566.73334
$ws.Range("I23").Value = 399.33334
$ws.Range("J23").Value = 608.5833
$ws.Range("K23").Value = 1198.00002
$ws.Range("L23").Value = 1825.7499
$ws.Range("M23").Value = -963.0000199999999
$ws.Range("N23").Value = -2295.7499

$ws.Range("H132").Value = 1718.6875
$ws.Range("I132").Value = 1528.4286
$ws.Range("K132").Value = 13755.8574
$ws.Range("M132").Value = -11225.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H97").Value = 2646519
$ws.Range("I97").Value = 4762574.5
$ws.Range("J97").Value = 1449.5
$ws.Range("K97").Value = 4762574.5
$ws.Range("L97").Value = 1449.5
$ws.Range("M97").Value = -4762078.5
$ws.Range("N97").Value = -2441.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7795594.5
$ws.Range("J2").Value = 28070.8
$ws.Range("L2").Value = 28070.8
$ws.Range("N2").Value = -28294.8

$ws.Range("H22").Value = 43502.477
$ws.Range("I22").Value = 89710.7
$ws.Range("J22").Value = 1495
$ws.Range("K22").Value = 89710.7
$ws.Range("L22").Value = 1495
$ws.Range("M22").Value = -89415.7
$ws.Range("N22").Value = -2085

$ws.Range("H27").Value = 43502.477
$ws.Range("I27").Value = 89710.7
$ws.Range("J27").Value = 1495
$ws.Range("K27").Value = 89710.7
$ws.Range("L27").Value = 1495
$ws.Range("M27").Value = -89603.7
$ws.Range("N27").Value = -1709

$ws.Range("H46").Value = 6035.3887
$ws.Range("I46").Value = 2050
$ws.Range("J46").Value = 6832.467
$ws.Range("K46").Value = 2050
$ws.Range("L46").Value = 6832.467
$ws.Range("M46").Value = -1862
$ws.Range("N46").Value = -7208.467

$ws.Range("H64").Value = 500010000
$ws.Range("J64").Value = 500010000
$ws.Range("L64").Value = 500010000
$ws.Range("N64").Value = -500010450

$ws.Range("H67").Value = 500010000
$ws.Range("J67").Value = 500010000
$ws.Range("L67").Value = 500010000
$ws.Range("N67").Value = -500011560

$ws.Range("H68").Value = 903
$ws.Range("J68").Value = 903
$ws.Range("L68").Value = 903
$ws.Range("N68").Value = -2401

$ws.Range("H71").Value = 903
$ws.Range("J71").Value = 903
$ws.Range("L71").Value = 4515
$ws.Range("N71").Value = -12003

$ws.Range("H75").Value = 24000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 24000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 24000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -25872

$ws.Range("H78").Value = 24000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 24000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 72000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -81360

$ws.Range("H100").Value = 3690.524
$ws.Range("I100").Value = 3307.923
$ws.Range("J100").Value = 4312.25
$ws.Range("K100").Value = 3307.923
$ws.Range("L100").Value = 4312.25
$ws.Range("M100").Value = -2766.923
$ws.Range("N100").Value = -5394.25

$ws.Range("H132").Value = 7533.7
$ws.Range("I132").Value = 7938.44
$ws.Range("K132").Value = 23815.32
$ws.Range("M132").Value = -21285.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9051.317999999999
$ws.Range("I62").Value = 4600
$ws.Range("J62").Value = 9496.450000000001
$ws.Range("K62").Value = 4600
$ws.Range("L62").Value = 9496.450000000001
$ws.Range("M62").Value = -3976
$ws.Range("N62").Value = -10744.45

$ws.Range("H65").Value = 9051.317999999999
$ws.Range("I65").Value = 4600
$ws.Range("J65").Value = 9496.450000000001
$ws.Range("K65").Value = 23000
$ws.Range("L65").Value = 47482.25
$ws.Range("M65").Value = -19880
$ws.Range("N65").Value = -53722.25

$ws.Range("H107").Value = 71430130
$ws.Range("I107").Value = 83333840
$ws.Range("J107").Value = 7822.5
$ws.Range("K107").Value = 250001520
$ws.Range("L107").Value = 23467.5
$ws.Range("M107").Value = -249999600
$ws.Range("N107").Value = -27307.5

$ws.Range("H136").Value = 3227.8333
$ws.Range("I136").Value = 2880.2058
$ws.Range("J136").Value = 4705.25
$ws.Range("K136").Value = 8640.617400000001
$ws.Range("L136").Value = 14115.75
$ws.Range("M136").Value = -6090.617400000001
$ws.Range("N136").Value = -19215.75
